$wb = $excel.ActiveWorkbook

$wsSurvey  = $wb.Worksheets.Item("eCL Survey")
$wsHistory = $wb.Worksheets.Item("Change History")

# --- "eCL Survey" sheet: update the "Quality Now" survey question row (row 8) ---
# Phase II site list replaces the single "London" site.
$wsSurvey.Range("B8").Value = "Chester, Hattiesburg, London, Lynn Haven, Tampa, Winchester,"
$wsSurvey.Range("B8").WrapText = $true

# --- "Change History" sheet: add new row 12 entry for the tfs20256 release ---
$wsHistory.Range("B12").Value = "3/31/2021"
$wsHistory.Range("C12").Value = "tfs20256 - ecl qn surverys - phase ii sites"
$wsHistory.Range("D12").Value = 1.07
$wsHistory.Range("E12").Value = "Doug Stearns"

# End date for the "Quality Now" question is no longer tracked as a date.
$wsSurvey.Range("M8").Value = "none"

# --- restore the selections shown in the saved workbook ---
[void]$wsHistory.Range("C13").Select()
[void]$wsSurvey.Range("A8").Select()
